$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 cells (field value edits per diff) ---
$ws.Range("G2").Value = '['''', ''Penang'', ''Bukit Mertajam'']'
$ws.Range("H2").Value = '[{''field_of_study'': ''Bachelor Of Computer Science (Data Engineering)'', ''level'': "Bachelor''s Degree", ''cgpa'': ''3.97'', ''university'': ''Universiti Teknologi Malaysia'', ''start_date'': ''2020'', ''year_of_graduation'': ''2024''}, {''field_of_study'': ''Foundation in Science'', ''level'': ''Foundation'', ''cgpa'': ''3.78'', ''university'': ''Universiti Teknologi Malaysia'', ''start_date'': ''2019'', ''year_of_graduation'': ''2020''}]'
$ws.Range("J2").Value = '[''Time Management'', ''Collaboration'', ''Adaptability'', ''Leadership'', ''Communication'']'
$ws.Range("K2").Value = '[''HTML 5'', ''CSS'', ''JavaScript'', ''PHP'', ''SQL'', ''Python'', ''.NET'', ''React'', ''spaCy'', ''NLTK'', ''TensorFlow'', ''PyTorch'', ''LangChain'', ''Llama'', ''Django'', ''PostgreSQL'', ''Laravel'', ''MySQL'', ''Microsoft SQL Server'', ''.NET MVC Framework'']'
$ws.Range("M2").Value = '[{''job_title'': ''Data Science Intern'', ''job_company'': ''Petronas Digital Sdn Bhd'', ''Industries'': ''Oil & Gas'', ''start_date'': ''2023-09'', ''end_date'': ''2024-06'', ''job_location'': ''Bukit Mertajam, Penang''}]'

# --- Add new row 3 (new candidate record) ---
# Copy A2 (number style s="1") to A3 so the new index cell matches formatting
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 'Ang Teik Hun'
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '+60124773683'
$ws.Range("D3").Value = 'teikhun0422@hotmail.com'
$ws.Range("E3").Value = 'Yes'
$ws.Range("F3").Value = 'N/A'
$ws.Range("G3").Value = '[{''Country'': ''Malaysia'', ''State'': ''Penang'', ''City'': ''Butterworth''}]'
$ws.Range("H3").Value = '[{''field_of_study'': ''Applied Data Analytics'', ''level'': "Master''s", ''cgpa'': ''6.42/7'', ''university'': ''Australian National University'', ''start_date'': ''2021-07'', ''year_of_graduation'': ''2022''}]'
$ws.Range("I3").Value = '[''Microsoft Certified: Azure Fundamentals'', ''Microsoft Certified: Azure AI Engineer Associate'']'
$ws.Range("J3").Value = '[''Analytical Thinking Skills'', ''Adaptability'', ''Time Management'', ''Leadership'', ''Power BI'', ''Python'', ''Neural Network'', ''Machine Learning (SKlearn)'', ''Database SQL'', ''Data Wrangling'', ''Optimization'']'
$ws.Range("K3").Value = '[''N/A'']'
$ws.Range("L3").Value = '[''Chinese'', ''English'', ''Malay'', ''French'']'
$ws.Range("M3").Value = '[{''job_title'': ''Data Scientist'', ''job_company'': ''Petroliam Nasional Berhad Group Digital'', ''Industries'': ''Oil & Gas'', ''start_date'': ''2020-11'', ''end_date'': ''2021-07'', ''job_location'': ''KL'', ''job_duration'': ''0.6''}]'
